# Adds a new "tpu py_file and data argumentation" block of results to the
# log sheet: one new header row (85) plus six new data rows (86-91),
# replacing the old trailing partial row that used to live at row 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear out the old tail row (previously r=86, a partial/leftover row
#    with only E:I populated) so nothing of it survives.
# ---------------------------------------------------------------------
$ws.Range("A86:O86").Clear()

# ---------------------------------------------------------------------
# 2. Row 85: new header row - copy formatting from the look-alike header
#    at row 76 (bold style) and set its own header captions.
# ---------------------------------------------------------------------
$ws.Range("B76:M76").Copy()
$ws.Range("B85:M85").PasteSpecial(-4122)

$ws.Range("B85").Value = "lr"
$ws.Range("C85").Value = "batch_size"
$ws.Range("D85").Value = "epoches"
$ws.Range("E85").Value = "epoch1"
$ws.Range("F85").Value = "epoch2"
$ws.Range("G85").Value = "epoch3"
$ws.Range("H85").Value = "epoch4"
$ws.Range("I85").Value = "epoch5"
$ws.Range("J85").Value = "mean"
$ws.Range("K85").Value = "focal_r"
$ws.Range("L85").Value = "label_smoothing"
$ws.Range("M85").Value = "score"

# ---------------------------------------------------------------------
# 3. Rows 86-91: new data rows. Copy formatting from row 77, an existing
#    data row that already has the right per-column styles (A wrap-text,
#    B percent/sci numfmt, L left aligned, rest default).
# ---------------------------------------------------------------------
$ws.Range("A77:M77").Copy()
$ws.Range("A86:M91").PasteSpecial(-4122)

$labels = @(
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.去除标点符号",
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.去除空格",
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.data argumentation",
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.data argumentation(使用了数据预处理)",
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.data argumentation(仅仅补充头部和尾部信息)",
  "1. roberta + 5-fold（使用别人的分词工具）+ conv1d（使用最后两层layer） 2.使用learing schedule 3.focal loss 4.使用logit比较 5.data argumentation"
)
$rowHeights = @(42, 42, 42, 56, 56, 42)

for ($i = 0; $i -lt 6; $i++) {
  $r = 86 + $i
  $ws.Range("A${r}").Value = $labels[$i]
  $ws.Range("J${r}").Formula = "=AVERAGE(E${r}:I${r})"
  $ws.Range("K${r}").Value = 1
  $ws.Range("L${r}").Value = "None"
  $ws.Rows.Item($r).RowHeight = $rowHeights[$i]
}

# Row 86
$ws.Range("B86").Value = 0.00004
$ws.Range("C86").Value = 64
$ws.Range("D86").Value = 5
$ws.Range("E86").Value = 0.64483199999999996
$ws.Range("F86").Value = 0.65651300000000001
$ws.Range("G86").Value = 0.65056899999999995
$ws.Range("H86:I86").Clear()
$ws.Range("M86").Value = 0.66300000000000003

# Row 87
$ws.Range("B87").Value = 0.00004
$ws.Range("C87").Value = 64
$ws.Range("D87").Value = 5
$ws.Range("E87").Value = 0.68095899999999998
$ws.Range("F87").Value = 0.69272599999999995
$ws.Range("G87").Value = 0.69733100000000003
$ws.Range("H87").Value = 0.69836799999999999
$ws.Range("I87").Value = 0.69530700000000001
$ws.Range("M87").Clear()

# Row 88
$ws.Range("B88").Value = 0.00004
$ws.Range("C88").Value = 64
$ws.Range("D88").Value = 5
$ws.Range("E88").Value = 0.697824
$ws.Range("F88").Value = 0.71689599999999998
$ws.Range("G88").Value = 0.72333800000000004
$ws.Range("H88").Value = 0.72910399999999997
$ws.Range("I88").Value = 0.72910399999999997
$ws.Range("M88").Clear()

# Row 89
$ws.Range("B89").Value = 0.00005
$ws.Range("C89").Value = 86
$ws.Range("D89").Value = 5
$ws.Range("E89").Value = 0.70022700000000004
$ws.Range("F89").Value = 0.72020099999999998
$ws.Range("G89").Value = 0.71912200000000004
$ws.Range("H89:I89").Clear()
$ws.Range("M89").Value = 0.71199999999999997

# Row 90
$ws.Range("B90").Value = 0.00004
$ws.Range("C90").Value = 96
$ws.Range("D90").Value = 5
$ws.Range("E90").Value = 0.69025599999999998
$ws.Range("F90").Value = 0.70514100000000002
$ws.Range("G90").Value = 0.71001999999999998
$ws.Range("H90").Value = 0.71414599999999995
$ws.Range("I90").Value = 0.71295600000000003
$ws.Range("M90").Value = 0.70799999999999996

# Row 91
$ws.Range("B91").Value = 0.00004
$ws.Range("C91").Value = 96
$ws.Range("D91").Value = 5
$ws.Range("E91").Value = 0.68267100000000003
$ws.Range("F91").Value = 0.69822499999999998
$ws.Range("G91").Value = 0.69930700000000001
$ws.Range("H91").Value = 0.70250299999999999
$ws.Range("I91").Value = 0.70250299999999999
$ws.Range("M91").Clear()

# ---------------------------------------------------------------------
# 4. Update the view: active selection to match the post-edit state
#    (user was looking near the newly appended rows).
# ---------------------------------------------------------------------
$ws.Range("J93").Select()
